$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D8").Value = 22688900
$ws.Range("E8").Value = 21609500
$ws.Range("F8").Value = 25275100
$ws.Range("G8").Value = 40294000
$ws.Range("H8").Value = 47588100
$ws.Range("I8").Value = 47852000
$ws.Range("J8").Value = 38207300
$ws.Range("D9").Value = 15680900
$ws.Range("E9").Value = 18585800
$ws.Range("F9").Value = 22389300
$ws.Range("G9").Value = 35976600
$ws.Range("H9").Value = 42366300
$ws.Range("I9").Value = 42333700
$ws.Range("J9").Value = 33491800
$ws.Range("D10").Value = 7007900
$ws.Range("E10").Value = 3023800
$ws.Range("F10").Value = 2885800
$ws.Range("G10").Value = 4317400
$ws.Range("H10").Value = 5221700
$ws.Range("I10").Value = 5518300
$ws.Range("J10").Value = 4715600
$ws.Range("D12").Value = 133500
$ws.Range("E12").Value = 199700
$ws.Range("F12").Value = 313000
$ws.Range("G12").Value = 253600
$ws.Range("H12").Value = 308500
$ws.Range("I12").Value = 327400
$ws.Range("J12").Value = 203200
$ws.Range("D14").Value = 1420400
$ws.Range("E14").Value = 2016200
$ws.Range("F14").Value = 3374900
$ws.Range("G14").Value = 1431700
$ws.Range("H14").Value = 539100
$ws.Range("I14").Value = 288600
$ws.Range("J14").Value = 420900
$ws.Range("D15").Value = 1976900
$ws.Range("D17").Value = 21316700
$ws.Range("E17").Value = 22122300
$ws.Range("F17").Value = 27513400
$ws.Range("G17").Value = 39207900
$ws.Range("H17").Value = 44733700
$ws.Range("I17").Value = 44369700
$ws.Range("J17").Value = 35408900
$ws.Range("D18").Value = 1372200
$ws.Range("E18").Value = -512700
$ws.Range("F18").Value = -2238400
$ws.Range("G18").Value = 1086100
$ws.Range("H18").Value = 2854300
$ws.Range("I18").Value = 3482300
$ws.Range("J18").Value = 2798400
$ws.Range("D20").Value = 434200
$ws.Range("E20").Value = 388200
$ws.Range("F20").Value = 437600
$ws.Range("G20").Value = 208700
$ws.Range("H20").Value = 57200
$ws.Range("I20").Value = 187900
$ws.Range("J20").Value = 122900
$ws.Range("D21").Value = 4042300
$ws.Range("E21").Value = 4131300
$ws.Range("F21").Value = 4001500
$ws.Range("G21").Value = 4856700
$ws.Range("H21").Value = 5493900
$ws.Range("I21").Value = 5959900
$ws.Range("J21").Value = "NA"
$ws.Range("D22").Value = 139100
$ws.Range("E22").Value = 133500
$ws.Range("F22").Value = 341100
$ws.Range("G22").Value = 406200
$ws.Range("H22").Value = 341100
$ws.Range("I22").Value = 464200
$ws.Range("J22").Value = 429500
$ws.Range("D23").Value = 1667300
$ws.Range("E23").Value = -258100
$ws.Range("F23").Value = -2141900
$ws.Range("G23").Value = 888600
$ws.Range("H23").Value = 2570500
$ws.Range("I23").Value = 3206100
$ws.Range("J23").Value = 2491800
$ws.Range("D24").Value = 711300
$ws.Range("E24").Value = -52700
$ws.Range("F24").Value = -733800
$ws.Range("G24").Value = 297300
$ws.Range("H24").Value = 630600
$ws.Range("I24").Value = 1197200
$ws.Range("J24").Value = 710200
$ws.Range("D26").Value = 955900
$ws.Range("E26").Value = -205300
$ws.Range("F26").Value = -1408100
$ws.Range("G26").Value = 591300
$ws.Range("H26").Value = 1939900
$ws.Range("I26").Value = 2008900
$ws.Range("J26").Value = 1781600
$ws.Range("D27").Value = 488100
$ws.Range("E27").Value = -452200
$ws.Range("F27").Value = -1234200
$ws.Range("G27").Value = 311900
$ws.Range("H27").Value = 1303800
$ws.Range("I27").Value = 1529700
$ws.Range("J27").Value = 1210600
$ws.Range("D32").Value = -434200
$ws.Range("E32").Value = -388200
$ws.Range("F32").Value = -437600
$ws.Range("G32").Value = -208700
$ws.Range("H32").Value = -57200
$ws.Range("I32").Value = -187900
$ws.Range("J32").Value = -122900
$ws.Range("D33").Value = 488100
$ws.Range("E33").Value = -452200
$ws.Range("F33").Value = -1234200
$ws.Range("G33").Value = 311900
$ws.Range("H33").Value = 1303800
$ws.Range("I33").Value = 1529700
$ws.Range("J33").Value = 1210600
$ws.Range("D35").Value = 488100
$ws.Range("E35").Value = -452200
$ws.Range("F35").Value = -1234200
$ws.Range("G35").Value = 311900
$ws.Range("H35").Value = 1303800
$ws.Range("I35").Value = 1529700
$ws.Range("J35").Value = 1210600
$ws.Range("D41").Value = 4456500
$ws.Range("E41").Value = 2321400
$ws.Range("F41").Value = 1512400
$ws.Range("G41").Value = 728200
$ws.Range("H41").Value = 791000
$ws.Range("I41").Value = 1377000
$ws.Range("J41").Value = 402600
$ws.Range("D42").Value = 20200
$ws.Range("E42").Value = 39300
$ws.Range("F42").Value = 68400
$ws.Range("G42").Value = 72900
$ws.Range("H42").Value = 936500
$ws.Range("I42").Value = 555800
$ws.Range("J42").Value = 430300
$ws.Range("D43").Value = 3287400
$ws.Range("E43").Value = 3330100
$ws.Range("F43").Value = 3617300
$ws.Range("G43").Value = 4373500
$ws.Range("H43").Value = 8160600
$ws.Range("I43").Value = 9469700
$ws.Range("J43").Value = 8302600
$ws.Range("D44").Value = 1686400
$ws.Range("E44").Value = 1865900
$ws.Range("F44").Value = 2101500
$ws.Range("G44").Value = 2503200
$ws.Range("H44").Value = 5510700
$ws.Range("I44").Value = 6811900
$ws.Range("J44").Value = 7066300
$ws.Range("D45").Value = 1093900
$ws.Range("E45").Value = 1044600
$ws.Range("F45").Value = 2255200
$ws.Range("G45").Value = 1632500
$ws.Range("H45").Value = 691600
$ws.Range("I45").Value = 944900
$ws.Range("J45").Value = 381000
$ws.Range("D46").Value = 10544500
$ws.Range("E46").Value = 8601200
$ws.Range("F46").Value = 9554900
$ws.Range("G46").Value = 9310300
$ws.Range("H46").Value = 8486700
$ws.Range("I46").Value = 10311900
$ws.Range("J46").Value = 8811200
$ws.Range("D47").Value = 5394500
$ws.Range("E47").Value = 4228800
$ws.Range("F47").Value = 3747400
$ws.Range("G47").Value = 3148300
$ws.Range("H47").Value = 3495700
$ws.Range("I47").Value = 5296300
$ws.Range("J47").Value = 3271200
$ws.Range("D48").Value = 15319700
$ws.Range("E48").Value = 16395600
$ws.Range("F48").Value = 18445500
$ws.Range("G48").Value = 20744500
$ws.Range("H48").Value = 38261800
$ws.Range("I48").Value = 18973500
$ws.Range("J48").Value = 31373500
$ws.Range("D49").Value = 2971000
$ws.Range("E49").Value = 1922000
$ws.Range("F49").Value = 3674500
$ws.Range("G49").Value = 3873100
$ws.Range("H49").Value = 8080200
$ws.Range("I49").Value = 7808100
$ws.Range("J49").Value = 7690400
$ws.Range("D52").Value = 1198300
$ws.Range("E52").Value = 4881800
$ws.Range("F52").Value = 1226300
$ws.Range("G52").Value = 907700
$ws.Range("H52").Value = 1423000
$ws.Range("I52").Value = 798300
$ws.Range("J52").Value = 396700
$ws.Range("D54").Value = 35428000
$ws.Range("E54").Value = 36029300
$ws.Range("F54").Value = 36648700
$ws.Range("G54").Value = 37985000
$ws.Range("H54").Value = 35733100
$ws.Range("I54").Value = 34242200
$ws.Range("J54").Value = 31879600
$ws.Range("D57").Value = 3659900
$ws.Range("E57").Value = 4186100
$ws.Range("F57").Value = 3792300
$ws.Range("G57").Value = 4858200
$ws.Range("H57").Value = 5513500
$ws.Range("I57").Value = 4813800
$ws.Range("J57").Value = 3849800
$ws.Range("D58").Value = 1028900
$ws.Range("E58").Value = 309700
$ws.Range("F58").Value = 568800
$ws.Range("G58").Value = 686700
$ws.Range("H58").Value = 1131000
$ws.Range("I58").Value = 886500
$ws.Range("J58").Value = 627800
$ws.Range("D59").Value = 2969900
$ws.Range("E59").Value = 3051800
$ws.Range("F59").Value = 4638300
$ws.Range("G59").Value = 4399300
$ws.Range("H59").Value = 4384200
$ws.Range("I59").Value = 5706000
$ws.Range("J59").Value = 4040600
$ws.Range("D60").Value = 7658700
$ws.Range("E60").Value = 7547600
$ws.Range("F60").Value = 8999500
$ws.Range("G60").Value = 9944200
$ws.Range("H60").Value = 9264300
$ws.Range("I60").Value = 8056100
$ws.Range("J60").Value = 7282400
$ws.Range("D61").Value = 5687400
$ws.Range("E61").Value = 5609900
$ws.Range("F61").Value = 5464100
$ws.Range("G61").Value = 5528000
$ws.Range("H61").Value = 4564300
$ws.Range("I61").Value = 5115900
$ws.Range("J61").Value = 4808300
$ws.Range("D62").Value = 5928600
$ws.Range("E62").Value = 6006000
$ws.Range("F62").Value = 6107000
$ws.Range("G62").Value = 6196800
$ws.Range("H62").Value = 5666700
$ws.Range("I62").Value = 5805800
$ws.Range("J62").Value = 4758400
$ws.Range("D66").Value = 22843700
$ws.Range("E66").Value = 23782800
$ws.Range("F66").Value = 23552800
$ws.Range("G66").Value = 24981100
$ws.Range("H66").Value = 22702300
$ws.Range("I66").Value = 20887800
$ws.Range("J66").Value = 19656300
$ws.Range("D72").Value = 9714200
$ws.Range("E72").Value = 9413500
$ws.Range("F72").Value = 10169700
$ws.Range("G72").Value = 11750600
$ws.Range("H72").Value = 23670400
$ws.Range("I72").Value = 23213300
$ws.Range("J72").Value = 22681900
$ws.Range("D76").Value = 12584200
$ws.Range("E76").Value = 12246500
$ws.Range("F76").Value = 13095900
$ws.Range("G76").Value = 13003900
$ws.Range("H76").Value = 13030800
$ws.Range("I76").Value = 13354400
$ws.Range("J76").Value = 12223300
$ws.Range("D81").Value = 488100
$ws.Range("E81").Value = -452200
$ws.Range("F81").Value = -1234200
$ws.Range("G81").Value = 311900
$ws.Range("H81").Value = 1303800
$ws.Range("I81").Value = 1529700
$ws.Range("J81").Value = 1210600
$ws.Range("D83").Value = 2230500
$ws.Range("E83").Value = 4245600
$ws.Range("F83").Value = 5788300
$ws.Range("G83").Value = 3553300
$ws.Range("H83").Value = 2576100
$ws.Range("I83").Value = 2284200
$ws.Range("J83").Value = "NA"
$ws.Range("D89").Value = 3868600
$ws.Range("E89").Value = 3229100
$ws.Range("F89").Value = 3180800
$ws.Range("G89").Value = 4113200
$ws.Range("H89").Value = 4627100
$ws.Range("I89").Value = 4278100
$ws.Range("J89").Value = 2820700
$ws.Range("D91").Value = -1779500
$ws.Range("E91").Value = -2268700
$ws.Range("F91").Value = -3341300
$ws.Range("G91").Value = -4301700
$ws.Range("H91").Value = -5349600
$ws.Range("I91").Value = -2788000
$ws.Range("J91").Value = -2762700
$ws.Range("D94").Value = -1981400
$ws.Range("E94").Value = -2016200
$ws.Range("F94").Value = -3224600
$ws.Range("G94").Value = -3808000
$ws.Range("H94").Value = -4466600
$ws.Range("I94").Value = -2557600
$ws.Range("J94").Value = "NA"
$ws.Range("D96").Value = -593500
$ws.Range("E96").Value = -522800
$ws.Range("F96").Value = -515000
$ws.Range("G96").Value = -513900
$ws.Range("H96").Value = -495900
$ws.Range("I96").Value = -702700
$ws.Range("J96").Value = -495300
$ws.Range("D100").Value = 30300
$ws.Range("E100").Value = -83000
$ws.Range("F100").Value = 850500
$ws.Range("G100").Value = -383700
$ws.Range("H100").Value = -719200
$ws.Range("I100").Value = -738200
$ws.Range("J100").Value = "NA"
$ws.Range("D101").Value = -47100
$ws.Range("E101").Value = -47100
$ws.Range("F101").Value = -21300
$ws.Range("G101").Value = 15700
$ws.Range("H101").Value = -26900
$ws.Range("I101").Value = -7900
$ws.Range("J101").Value = "NA"
$ws.Range("D102").Value = 1870400
$ws.Range("E102").Value = 1082700
$ws.Range("F102").Value = 785400
$ws.Range("G102").Value = -62800
$ws.Range("H102").Value = -585700
$ws.Range("I102").Value = 974400
$ws.Range("J102").Value = -658900
